# Updates the cryptos list figures (Price / Volume(1h)) to the latest
# scraped values, and fixes the Mantle/VeChain row ordering.
#
# Notes:
#  - Column D (Price) and E (Volume 1h) are stored as plain text in the
#    sheet (e.g. "54.129.27", "  -10.79%  "), not numbers. Several of the
#    new Price values (e.g. "0.999", "10.28") look like valid numbers, so
#    a plain assignment would make Excel silently convert them to numeric
#    cells. Prefixing those values with a leading apostrophe forces Excel
#    to keep (and store) them as text, matching the original formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    if ($value -match '^[+-]?[0-9]*\.?[0-9]+$') {
        # Value looks like a plain number (e.g. "0.999"); prefix with an
        # apostrophe so Excel keeps storing it as text, same as the
        # original cell, instead of silently converting it to a number.
        $ws.Range($range).Value = "'" + $value
    } else {
        $ws.Range($range).Value = $value
    }
}

# row -> (Price, Volume(1h))
$updates = @{
    2  = @('54.177.98',  '  -10.67%  ')
    3  = @('2.309.17',   '  -20.46%  ')
    4  = @('0.999',      '  -0.11%  ')
    5  = @('448.20',     '  -15.15%  ')
    6  = @('128.86',     '  -10.68%  ')
    7  = @('0.996',      '  -0.29%  ')
    8  = @('0.475',      '  -14.51%  ')
    9  = @('2.280.66',   '  -21.63%  ')
    10 = @('5.39',       '  -10.53%  ')
    11 = @('0.0924',     '  -14.91%  ')
    12 = @('0.311',      '  -14.32%  ')
    13 = @($null,        '  -3.53%  ')
    14 = @('2.697.75',   '  -20.87%  ')
    15 = @('54.179.87',  '  -10.64%  ')
    16 = @('18.82',      '  -17.40%  ')
    17 = @('0.0000121',  '  -14.31%  ')
    18 = @('2.318.42',   '  -20.26%  ')
    19 = @('4.06',       '  -19.56%  ')
    20 = @('9.44',       '  -19.41%  ')
    21 = @('300.08',     '  -17.07%  ')
    22 = @($null,        '  -0.09%  ')
    23 = @('5.61',       '  -1.31%  ')
    24 = @('5.34',       '  -19.46%  ')
    25 = @('55.70',      '  -14.09%  ')
    26 = @('0.985',      '  -1.33%  ')
    27 = @($null,        '  -13.76%  ')
    28 = @('0.371',      '  -18.80%  ')
    29 = @('0.996',      '  -0.33%  ')
    30 = @('6.81',       '  -13.67%  ')
    31 = @('0.0₃0711',   '  -17.47%  ')
    32 = @('146.61',     '  -3.22%  ')
    33 = @('16.96',      '  -14.05%  ')
    34 = @($null,        '  -19.40%  ')
    35 = @($null,        '  -15.63%  ')
    36 = @('3.62',       '  -17.72%  ')
    37 = @('0.843',      '  -16.73%  ')
    38 = @($null,        '  -16.19%  ')
    39 = @('0.993',      '  -0.45%  ')
    40 = @('32.99',      '  -12.44%  ')
    41 = @('10.28',      '  +0.40%  ')
    42 = @($null,        '  -17.24%  ')
    43 = @('3.15',       '  -15.57%  ')
    44 = @('1.930.63',   '  -15.72%  ')
    45 = @('0.0497',     '  -14.74%  ')
    48 = @('0.0824',     '  -11.01%  ')
    49 = @('16.10',      '  -21.78%  ')
    50 = @('4.02',       '  -20.01%  ')
    51 = @('4.69',       '  -3.01%  ')
}

foreach ($row in $updates.Keys | Sort-Object) {
    $pair = $updates[$row]
    $price = $pair[0]
    $volume = $pair[1]
    if ($null -ne $price) {
        Set-TextValue "D$row" $price
    }
    Set-TextValue "E$row" $volume
}

# Rows 46/47 were swapped (Mantle <-> VeChain) and refreshed with new
# price/volume figures.
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D46' '0.0209'
Set-TextValue 'E46' '  -12.33%  '

$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D47' '0.515'
Set-TextValue 'E47' '  -20.35%  '
